{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Changes applied (per the commit's XML diff):\n//   1. Demote the three Heading2 paragraphs (\"Introduction\",\n//      \"Critical Challenges in Rural Healthcare\", \"Conclusion\") back to\n//      the default/\"Normal\" paragraph style (their <w:pPr>/<w:pStyle>\n//      is dropped in the target OOXML).\n//   2. Replace in-text citation markers with the updated reference tags:\n//        (Wang)            -> (Ref-u378155)\n//        (Maganty et al.)  -> (Ref-u378155) [1st occurrence]\n//                          -> (Nguyen, 2015) [2nd occurrence]\n//        (Coombs et al.)   -> (Nguyen, 2015) [both occurrences]\n//        (Palozzi et al.)  -> (Ref-f737193)\n//        (Hirko et al.)    -> (Ref-f737193)\n\nconst body = context.document.body;\n\n// --- 1. Strip the Heading2 style from the three section headings. ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst headingTexts = [\n  \"Introduction\",\n  \"Critical Challenges in Rural Healthcare\",\n  \"Conclusion\",\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const trimmed = para.text.trim();\n  if (headingTexts.indexOf(trimmed) !== -1) {\n    para.style = \"Normal\";\n  }\n}\nawait context.sync();\n\n// --- 2. Update the citation markers. ---\n// Run every search against the (still untouched) body text first so the\n// returned ranges are all resolved against stable, original offsets; only\n// afterwards do we mutate them. This avoids occurrence-index drift that\n// would happen if we re-searched after each edit.\nconst searchWang = body.search(\"(Wang)\", { matchCase: true });\nconst searchMaganty = body.search(\"(Maganty et al.)\", { matchCase: true });\nconst searchCoombs = body.search(\"(Coombs et al.)\", { matchCase: true });\nconst searchPalozzi = body.search(\"(Palozzi et al.)\", { matchCase: true });\nconst searchHirko = body.search(\"(Hirko et al.)\", { matchCase: true });\n\nsearchWang.load(\"text\");\nsearchMaganty.load(\"text\");\nsearchCoombs.load(\"text\");\nsearchPalozzi.load(\"text\");\nsearchHirko.load(\"text\");\nawait context.sync();\n\n// (Wang) -> (Ref-u378155)\nsearchWang.items[0].insertText(\"(Ref-u378155)\", \"Replace\");\n\n// (Maganty et al.) appears twice: 1st -> (Ref-u378155), 2nd -> (Nguyen, 2015)\nsearchMaganty.items[0].insertText(\"(Ref-u378155)\", \"Replace\");\nsearchMaganty.items[1].insertText(\"(Nguyen, 2015)\", \"Replace\");\n\n// (Coombs et al.) appears twice, both -> (Nguyen, 2015)\nsearchCoombs.items[0].insertText(\"(Nguyen, 2015)\", \"Replace\");\nsearchCoombs.items[1].insertText(\"(Nguyen, 2015)\", \"Replace\");\n\n// (Palozzi et al.) -> (Ref-f737193)\nsearchPalozzi.items[0].insertText(\"(Ref-f737193)\", \"Replace\");\n\n// (Hirko et al.) -> (Ref-f737193)\nsearchHirko.items[0].insertText(\"(Ref-f737193)\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $d = $word.ActiveDocument is pre-seeded by the harness.\n#\n# Changes applied (per the commit's XML diff):\n#   1. Demote the three Heading2 paragraphs (\"Introduction\",\n#      \"Critical Challenges in Rural Healthcare\", \"Conclusion\") back to\n#      the default/\"Normal\" paragraph style (their <w:pPr>/<w:pStyle>\n#      is dropped in the target OOXML).\n#   2. Replace in-text citation markers with the updated reference tags:\n#        (Wang)            -> (Ref-u378155)\n#        (Maganty et al.)  -> (Ref-u378155) [1st occurrence]\n#                          -> (Nguyen, 2015) [2nd occurrence]\n#        (Coombs et al.)   -> (Nguyen, 2015) [both occurrences]\n#        (Palozzi et al.)  -> (Ref-f737193)\n#        (Hirko et al.)    -> (Ref-f737193)\n\n$d = $word.ActiveDocument\n\n# --- 1. Strip the Heading2 style from the three section headings. ---\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.Trim()\n    if ($text -eq \"Introduction\" -or $text -eq \"Critical Challenges in Rural Healthcare\" -or $text -eq \"Conclusion\") {\n        $p.Style = \"Normal\"\n    }\n}\n\n# --- 2. Update the citation markers. ---\n# wdReplaceNone = 0, wdReplaceOne = 1, wdReplaceAll = 2\n$wdReplaceAll = 2\n\n# Locate the two body paragraphs that hold the \"Maganty et al.\" /\n# \"Coombs et al.\" markers so each occurrence can be retargeted to its own\n# replacement text without clobbering the other mentions in the same\n# document.\n$paraRuralChallenges = $null\n$paraSocioeconomic = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if ($text -like \"*profound shortage of medical professionals*\") {\n        $paraRuralChallenges = $p\n    }\n    if ($text -like \"*inadequacy of healthcare access in rural areas is exacerbated*\") {\n        $paraSocioeconomic = $p\n    }\n}\n\n# (Wang) -> (Ref-u378155)  [whole-document: single occurrence]\n$range = $d.Content\n$range.Find.Execute(\"(Wang)\", $true, $false, $false, $false, $false, $true, 0, $false, \"(Ref-u378155)\", $wdReplaceAll)\n\n# (Maganty et al.), 1st occurrence (in the \"Critical Challenges\" paragraph) -> (Ref-u378155)\n$range = $paraRuralChallenges.Range\n$range.Find.Execute(\"(Maganty et al.)\", $true, $false, $false, $false, $false, $true, 0, $false, \"(Ref-u378155)\", $wdReplaceAll)\n\n# (Maganty et al.), 2nd occurrence (in the socioeconomic paragraph) -> (Nguyen, 2015)\n$range = $paraSocioeconomic.Range\n$range.Find.Execute(\"(Maganty et al.)\", $true, $false, $false, $false, $false, $true, 0, $false, \"(Nguyen, 2015)\", $wdReplaceAll)\n\n# (Coombs et al.), both occurrences (in the socioeconomic paragraph) -> (Nguyen, 2015)\n$range = $paraSocioeconomic.Range\n$range.Find.Execute(\"(Coombs et al.)\", $true, $false, $false, $false, $false, $true, 0, $false, \"(Nguyen, 2015)\", $wdReplaceAll)\n\n# (Palozzi et al.) -> (Ref-f737193)  [whole-document: single occurrence]\n$range = $d.Content\n$range.Find.Execute(\"(Palozzi et al.)\", $true, $false, $false, $false, $false, $true, 0, $false, \"(Ref-f737193)\", $wdReplaceAll)\n\n# (Hirko et al.) -> (Ref-f737193)  [whole-document: single occurrence]\n$range = $d.Content\n$range.Find.Execute(\"(Hirko et al.)\", $true, $false, $false, $false, $false, $true, 0, $false, \"(Ref-f737193)\", $wdReplaceAll)\n"}
